# Update "paises.xlsx" (COVID country stats) to the next scrape snapshot.
# - refresh the "last updated" timestamp
# - refresh case/death/recovery counters for a set of countries
# - two countries (Estado de Palestina, Angola) overtake their neighbours in
#   the ranking (sorted by "Casos totales" desc) and move up a few rows
# - two tied pairs (Fiyi/Dominica, Groenlandia/Islas Malvinas) swap places

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header timestamp
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "Datos actualizados a 29 de Junio de 2020 a las 20:51"

# ---------------------------------------------------------------------------
# 2) Helper to write a full data row (B..H) in one shot
# ---------------------------------------------------------------------------
function Set-Stats($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value2 = $b
    $ws.Cells.Item($row, 3).Value2 = $c
    $ws.Cells.Item($row, 4).Value2 = $d
    $ws.Cells.Item($row, 5).Value2 = $e
    $ws.Cells.Item($row, 6).Value2 = $f
    $ws.Cells.Item($row, 7).Value2 = $g
    $ws.Cells.Item($row, 8).Value2 = $h
}

# ---------------------------------------------------------------------------
# 3) In-place statistic refreshes (country keeps its row / rank)
# ---------------------------------------------------------------------------
Set-Stats 4   2653539 16462 1099297 1425674 0 131 128568   # Estados Unidos
Set-Stats 7   567233  18036 335272  215057  0 417 16904    # India
Set-Stats 17  195235  371   178100  8099    0 7   9036     # Alemania
Set-Stats 34  48246   449   37076   10856   0 1   314      # Emiratos Arabes Unidos
Set-Stats 74  8222    274   5496    2703    0 1   23       # Uzbekistan
Set-Stats 93  4237    88    1497    2612    0 2   128      # Mauritania
Set-Stats 102 2560    52    2325    203     0 0   32       # Mayotte
Set-Stats 121 1568    11    1311    235     0 0   22       # Zambia
Set-Stats 128 1152    6     260     879     0 0   13       # Malaui
Set-Stats 134 996     2     833     144     0 0   19       # Republica de Chipre
Set-Stats 140 866     0     781     11      0 0   74       # Republica del Chad
Set-Stats 143 795     14    380     404     0 0   11       # Suazilandia

# ---------------------------------------------------------------------------
# 4) Estado de Palestina jumps from row 110 up to row 106 (ahead of
#    Nicaragua, Mali, Madagascar and Paraguay, which each drop one place).
#    Insert a fresh row at the new position, fill it with Palestina's
#    updated numbers, then remove the row it vacated (now shifted to 111).
# ---------------------------------------------------------------------------
$ws.Rows("106:106").Insert()
$ws.Range("A106").Value2 = "Estado de Palestina"
Set-Stats 106 2185 195 447 1733 0 1 5
$ws.Rows("111:111").Delete()

# ---------------------------------------------------------------------------
# 5) Angola jumps from row 164 up to row 162 (ahead of Comoras and Siria,
#    which each drop one place). Same insert/fill/delete technique.
# ---------------------------------------------------------------------------
$ws.Rows("162:162").Insert()
$ws.Range("A162").Value2 = "Angola"
Set-Stats 162 276 9 93 172 0 0 11
$ws.Rows("165:165").Delete()

# ---------------------------------------------------------------------------
# 6) Tied pairs swap adjacent rows (no numeric change, just reordering).
# ---------------------------------------------------------------------------
$ws.Range("A205").Value2 = "Fiyi"
$ws.Range("A206").Value2 = "Dominica"

$ws.Range("A209").Value2 = "Groenlandia"
$ws.Range("A210").Value2 = "Islas Malvinas"
